$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# --- Row 14: height 15 -> 30 ; E14 text change ---
$ws.Range("A14").EntireRow.RowHeight = 30
$ws.Range("E14").Value = "El sistema generará un listado de las 10 placas mas buscada con los datos principales del vehiculo."

# --- Row 15: height 30 -> 15 ; C15, E15 text change ---
$ws.Range("A15").EntireRow.RowHeight = 15
$ws.Range("C15").Value = "Sistema"
$ws.Range("E15").Value = "el sistema permitira cargar los datos de una infraccion y registrarlos"

# --- Row 16: height 15 -> 30 ; C16, E16 text change ---
$ws.Range("A16").EntireRow.RowHeight = 30
$ws.Range("C16").Value = "Usuario solicitante"
$ws.Range("E16").Value = "El sistema permitira subir archivos multimedia como imágenes, audio y video."

# --- Row 17: C17 cleared ; E17 text change ---
$ws.Range("C17").Value = ""
$ws.Range("E17").Value = "El sistema debe permitir buscar por codigo de placa"

# --- Row 18: E18 text change only ---
$ws.Range("E18").Value = "El sistema debe permitir buscar por DNI de propietario"

# --- Row 19: height 30 -> 15 ; C19, E19 text change ---
$ws.Range("A19").EntireRow.RowHeight = 15
$ws.Range("C19").Value = "Sistema"
$ws.Range("E19").Value = "El sistema permitra realizar el pago de una infraccion"

# --- Row 20: height 15 -> 30 ; C20, E20 text change ---
$ws.Range("A20").EntireRow.RowHeight = 30
$ws.Range("C20").Value = "Usuario solicitante"
$ws.Range("E20").Value = "El sistema permitia el pago en linea a traves de pasarela de pagos de visa o mastercard."

# --- Row 22: new C22 and E22 values (copy style from row 21, BEFORE row 21 restyles) ---
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "Sistema"

$ws.Range("E21").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "El sistema debe permitir buscar por DNI de propietario"

# --- Row 21: E21 text + style change (s=22 -> s=12) ---
$ws.Range("E20").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "El sistema debe permitir buscar por codigo de placa"

# --- Row 26: remove F26:H26, add blank E26 ---
$ws.Range("F26:H26").Clear()
$ws.Range("E25").Copy()
$ws.Range("E26").PasteSpecial(-4122)

# --- sheetView: drop topLeftCell, move selection to E26 ---
$ws.Range("E26").Select()

Write-Host "edit complete"
